$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.128.44"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "3.369.76"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.54"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.94"
$ws.Range("E6").Value = "  +9.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.592"
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.671"
$ws.Range("E9").Value = "  +3.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.120"
$ws.Range("E10").Value = "  +5.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.38"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("E12").Value = "  -0.81%  "

$ws.Range("D13").Value = "3.913.43"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.33"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.65"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "3.424.46"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "61.165.76"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  +1.25%  "

$ws.Range("E20").Value = "  +6.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.22"
$ws.Range("E21").Value = "  -5.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "83.50"
$ws.Range("E22").Value = "  +9.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.76"
$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "304.26"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.77"
$ws.Range("E26").Value = "  +11.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("E27").Value = "  +8.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.36"
$ws.Range("E28").Value = "  -5.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  -8.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.173"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("E31").Value = "  +1.58%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.32"
$ws.Range("E33").Value = "  -1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.54"
$ws.Range("E34").Value = "  -3.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0478"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.81"
$ws.Range("E37").Value = "  -0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -3.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -3.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.57"
$ws.Range("E41").Value = "  +2.11%  "

$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.288"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.98"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("E46").Value = "  -4.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  +2.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.38"
$ws.Range("E48").Value = "  -3.88%  "

$ws.Range("D49").Value = "2.117.30"
$ws.Range("E49").Value = "  -3.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  -4.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.89"
$ws.Range("E51").Value = "  +0.37%  "
